$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rows 158-199: the weekly Vega Monumental Concepcion / Betarraga
# series shifts down by one week (a new week is inserted at the top,
# rows 158-159). Only cells whose value actually differs from what
# was already sitting in that row are rewritten here (mostly Fecha,
# plus a few Volumen/Precio promedio/Precio $/Kg values).
$ws.Range("D158").Value = 44463
$ws.Range("O158").Value = 'Región de Coquimbo'
$ws.Range("D159").Value = 44463
$ws.Range("O159").Value = 'Región de Coquimbo'
$ws.Range("D160").Value = 44196
$ws.Range("D161").Value = 44196
$ws.Range("D162").Value = 44243
$ws.Range("D163").Value = 44243
$ws.Range("D164").Value = 44252
$ws.Range("J164").Value = 800
$ws.Range("D165").Value = 44252
$ws.Range("J165").Value = 400
$ws.Range("D166").Value = 44166
$ws.Range("D167").Value = 44166
$ws.Range("D168").Value = 44168
$ws.Range("J168").Value = 600
$ws.Range("D169").Value = 44168
$ws.Range("J169").Value = 300
$ws.Range("D170").Value = 44316
$ws.Range("J170").Value = 1000
$ws.Range("D171").Value = 44316
$ws.Range("J171").Value = 500
$ws.Range("D172").Value = 44397
$ws.Range("D173").Value = 44397
$ws.Range("D174").Value = 44273
$ws.Range("D175").Value = 44273
$ws.Range("D176").Value = 44372
$ws.Range("D177").Value = 44372
$ws.Range("D178").Value = 44365
$ws.Range("D179").Value = 44365
$ws.Range("D180").Value = 44306
$ws.Range("J180").Value = 600
$ws.Range("D181").Value = 44306
$ws.Range("J181").Value = 300
$ws.Range("D182").Value = 44215
$ws.Range("D183").Value = 44215
$ws.Range("D184").Value = 44357
$ws.Range("D185").Value = 44357
$ws.Range("D186").Value = 44162
$ws.Range("J186").Value = 800
$ws.Range("M186").Value = 650
$ws.Range("P186").Value = 130
$ws.Range("D187").Value = 44162
$ws.Range("J187").Value = 400
$ws.Range("D188").Value = 44239
$ws.Range("J188").Value = 700
$ws.Range("M188").Value = 643
$ws.Range("P188").Value = 129
$ws.Range("D189").Value = 44239
$ws.Range("D190").Value = 44376
$ws.Range("D191").Value = 44376
$ws.Range("D192").Value = 44292
$ws.Range("D193").Value = 44292
$ws.Range("D194").Value = 44358
$ws.Range("D195").Value = 44358
$ws.Range("D196").Value = 44211
$ws.Range("D197").Value = 44211
$ws.Range("D198").Value = 44425
$ws.Range("D199").Value = 44425

# --- New rows 200-201: Primera/Segunda calidad pair that used to be
# the last two rows (198-199) is appended again at the bottom with
# its original date/values, growing the sheet from R199 to R201.
$ws.Range("A200").Value = 11
$ws.Range("B200").Value = 'Vega Monumental Concepción'
$ws.Range("C200").Value = 'Bíobío'
$ws.Range("D200").Value = 44323
$ws.Range("E200").Value = 8
$ws.Range("F200").Value = 100114014
$ws.Range("G200").Value = 'Betarraga'
$ws.Range("H200").Value = 'Sin especificar'
$ws.Range("I200").Value = 'Primera'
$ws.Range("J200").Value = 600
$ws.Range("K200").Value = 600
$ws.Range("L200").Value = 700
$ws.Range("M200").Value = 650
$ws.Range("N200").Value = '$/paquete 5 unidades'
$ws.Range("O200").Value = 'Región Metropolitana'
$ws.Range("P200").Value = 130
$ws.Range("Q200").Value = 5
$ws.Range("R200").Value = 'Hortaliza'
$ws.Range("D200").NumberFormat = "YYYY-MM-DD HH:MM:SS"

$ws.Range("A201").Value = 11
$ws.Range("B201").Value = 'Vega Monumental Concepción'
$ws.Range("C201").Value = 'Bíobío'
$ws.Range("D201").Value = 44323
$ws.Range("E201").Value = 8
$ws.Range("F201").Value = 100114014
$ws.Range("G201").Value = 'Betarraga'
$ws.Range("H201").Value = 'Sin especificar'
$ws.Range("I201").Value = 'Segunda'
$ws.Range("J201").Value = 300
$ws.Range("K201").Value = 500
$ws.Range("L201").Value = 500
$ws.Range("M201").Value = 500
$ws.Range("N201").Value = '$/paquete 5 unidades'
$ws.Range("O201").Value = 'Región Metropolitana'
$ws.Range("P201").Value = 100
$ws.Range("Q201").Value = 5
$ws.Range("R201").Value = 'Hortaliza'
$ws.Range("D201").NumberFormat = "YYYY-MM-DD HH:MM:SS"

